# Edit "Analysis results.xlsx" Summary sheet: add a 4th scenario ("Only
# column creation" seems pre-existing; here we add the Hello-World /
# mvn exec:exec timing example), restructure the metrics table (drop the
# "Renaming columns, creating cache" row, reorder Clean-up before
# Transformations, add a per-row average column H), and make "Summary"
# the active/selected sheet again.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)                     # "Summary"
$ws2 = $wb.Worksheets.Item("Results per # of dataset and op")

# ---------------------------------------------------------------------
# 1. Clear out the old data table (rows 2-11, cols A-F) on the Summary
#    sheet, keeping row 1's headers untouched.
# ---------------------------------------------------------------------
$ws1.Range("A2:H13").Clear()

# ---------------------------------------------------------------------
# 2. Row 3: the three scenario sub-headers (mvn commands).
# ---------------------------------------------------------------------
$ws1.Range("A3").Value = "mvn exec:exec"

$ws1.Range("C3").Value = "mvn exec:exec -DexecMode=COL"
$ws1.Range("C3").Font.Name = "Monaco"
$ws1.Range("C3").Font.Size = 10

$ws1.Range("E3").Value = "mvn exec:exec -DexecMode=FULL"
$ws1.Range("E3").Font.Name = "Monaco"
$ws1.Range("E3").Font.Size = 10

# ---------------------------------------------------------------------
# 3. Rows 4-7: the four timed steps, each with its label + extraction
#    formula in B/D/F, and an H-column average of the three scenarios.
# ---------------------------------------------------------------------
$labels = @(
  @(4, "1. Creating a session .............. 1791", "1. Creating a session .............. 1553", "1. Creating a session .............. 1903"),
  @(5, "2. Loading initial dataset ......... 3287", "2. Loading initial dataset ......... 3197", "2. Loading initial dataset ......... 3184"),
  @(6, "3. Building full dataset ........... 242",  "3. Building full dataset ........... 208",  "3. Building full dataset ........... 213"),
  @(7, "4. Clean-up ........................ 8",     "4. Clean-up ........................ 8",     "4. Clean-up ........................ 8")
)

foreach ($row in $labels) {
  $r = $row[0]
  $ws1.Range("A$r").Value = $row[1]
  $ws1.Range("C$r").Value = $row[2]
  $ws1.Range("E$r").Value = $row[3]

  $ws1.Range("A$r").Font.Name = "Monaco"
  $ws1.Range("A$r").Font.Size = 10
  $ws1.Range("C$r").Font.Name = "Monaco"
  $ws1.Range("C$r").Font.Size = 10
  $ws1.Range("E$r").Font.Name = "Monaco"
  $ws1.Range("E$r").Font.Size = 10

  $ws1.Range("B$r").Formula = "=INT(RIGHT(A$r,LEN(A$r)-FIND(`"@`",SUBSTITUTE(A$r,`".`",`"@`",LEN(A$r)-LEN(SUBSTITUTE(A$r,`".`",`"`"))),1)-1))"
  $ws1.Range("D$r").Formula = "=INT(RIGHT(C$r,LEN(C$r)-FIND(`"@`",SUBSTITUTE(C$r,`".`",`"@`",LEN(C$r)-LEN(SUBSTITUTE(C$r,`".`",`"`"))),1)-1))"
  $ws1.Range("F$r").Formula = "=INT(RIGHT(E$r,LEN(E$r)-FIND(`"@`",SUBSTITUTE(E$r,`".`",`"@`",LEN(E$r)-LEN(SUBSTITUTE(E$r,`".`",`"`"))),1)-1))"

  $ws1.Range("H$r").Formula = "=AVERAGE(F$r,D$r,B$r)"
}
$ws1.Range("H5:H8").NumberFormat = "0"

# Row 8: blank label cells (still Monaco-styled, like the rest of the
# column) plus the running total of the averages so far.
$ws1.Range("A8").Font.Name = "Monaco"
$ws1.Range("A8").Font.Size = 10
$ws1.Range("C8").Font.Name = "Monaco"
$ws1.Range("C8").Font.Size = 10
$ws1.Range("E8").Font.Name = "Monaco"
$ws1.Range("E8").Font.Size = 10
$ws1.Range("H8").Formula = "=SUM(H4:H7)"

# ---------------------------------------------------------------------
# 4. Rows 9-10: Transformations + Final action.
# ---------------------------------------------------------------------
$ws1.Range("A9").Value = "5. Transformations`u{00A0} ................ 0"
$ws1.Range("C9").Value = "5. Transformations`u{00A0} ................ 182"
$ws1.Range("E9").Value = "5. Transformations`u{00A0} ................ 205"

$ws1.Range("A10").Value = "6. Final action .................... 20770"
$ws1.Range("C10").Value = "6. Final action .................... 34061"
$ws1.Range("E10").Value = "6. Final action .................... 24909"

foreach ($r in 9,10) {
  $ws1.Range("A$r").Font.Name = "Monaco"
  $ws1.Range("A$r").Font.Size = 10
  $ws1.Range("C$r").Font.Name = "Monaco"
  $ws1.Range("C$r").Font.Size = 10
  $ws1.Range("E$r").Font.Name = "Monaco"
  $ws1.Range("E$r").Font.Size = 10

  $ws1.Range("B$r").Formula = "=INT(RIGHT(A$r,LEN(A$r)-FIND(`"@`",SUBSTITUTE(A$r,`".`",`"@`",LEN(A$r)-LEN(SUBSTITUTE(A$r,`".`",`"`"))),1)-1))"
  $ws1.Range("D$r").Formula = "=INT(RIGHT(C$r,LEN(C$r)-FIND(`"@`",SUBSTITUTE(C$r,`".`",`"@`",LEN(C$r)-LEN(SUBSTITUTE(C$r,`".`",`"`"))),1)-1))"
  $ws1.Range("F$r").Formula = "=INT(RIGHT(E$r,LEN(E$r)-FIND(`"@`",SUBSTITUTE(E$r,`".`",`"@`",LEN(E$r)-LEN(SUBSTITUTE(E$r,`".`",`"`"))),1)-1))"
}

# ---------------------------------------------------------------------
# 5. Row 11: "Total processing time" label + SUM(row6:row10) per column.
# ---------------------------------------------------------------------
$ws1.Range("A11").Value = "Total processing time (excluding loading)"
$ws1.Range("B11").Formula = "=SUM(B6:B10)"
$ws1.Range("D11").Formula = "=SUM(D6:D10)"
$ws1.Range("F11").Formula = "=SUM(F6:F10)"
$ws1.Range("B11").Font.Bold = $true
$ws1.Range("D11").Font.Bold = $true
$ws1.Range("F11").Font.Bold = $true

# ---------------------------------------------------------------------
# 6. Row 12: "# of records" + extraction formula.
# ---------------------------------------------------------------------
$ws1.Range("A12").Value = "# of records ....................... 2487641"
$ws1.Range("C12").Value = "# of records ....................... 2487641"
$ws1.Range("E12").Value = "# of records ....................... 2487641"
$ws1.Range("A12").Font.Name = "Monaco"
$ws1.Range("A12").Font.Size = 10
$ws1.Range("C12").Font.Name = "Monaco"
$ws1.Range("C12").Font.Size = 10
$ws1.Range("E12").Font.Name = "Monaco"
$ws1.Range("E12").Font.Size = 10
$ws1.Range("B12").Formula = "=INT(RIGHT(A12,LEN(A12)-FIND(`"@`",SUBSTITUTE(A12,`".`",`"@`",LEN(A12)-LEN(SUBSTITUTE(A12,`".`",`"`"))),1)-1))"
$ws1.Range("D12").Formula = "=INT(RIGHT(C12,LEN(C12)-FIND(`"@`",SUBSTITUTE(C12,`".`",`"@`",LEN(C12)-LEN(SUBSTITUTE(C12,`".`",`"`"))),1)-1))"
$ws1.Range("F12").Formula = "=INT(RIGHT(E12,LEN(E12)-FIND(`"@`",SUBSTITUTE(E12,`".`",`"@`",LEN(E12)-LEN(SUBSTITUTE(E12,`".`",`"`"))),1)-1))"

# ---------------------------------------------------------------------
# 7. Row 13: ratio (# records / total processing time).
# ---------------------------------------------------------------------
$ws1.Range("B13").Formula = "=B12/B11"
$ws1.Range("D13").Formula = "=D12/D11"
$ws1.Range("F13").Formula = "=F12/F11"

# ---------------------------------------------------------------------
# 8. Column H width (new column, auto-fit like the diff's bestFit).
# ---------------------------------------------------------------------
$ws1.Range("H1").ColumnWidth = 12.66

# ---------------------------------------------------------------------
# 9. Selection + active sheet: "Summary" becomes the active tab again,
#    with F12 selected; "Results per # of dataset and op" loses focus.
# ---------------------------------------------------------------------
$ws1.Activate()
$ws1.Range("F12").Select()
